$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.339.42"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "'1.934.22"
$ws.Range("E3").Value = "  +0.14%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'0.7561"
$ws.Range("E5").Value = "  +5.29%  "
$ws.Range("D6").Value = "'245.06"
$ws.Range("E6").Value = "  -2.57%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").Value = "'0.3180"
$ws.Range("E8").Value = "  -2.71%  "
$ws.Range("E9").Value = "  +0.30%  "
$ws.Range("D10").Value = "'0.06990"
$ws.Range("E10").Value = "  -2.69%  "
$ws.Range("D11").Value = "'0.7774"
$ws.Range("E11").Value = "  -2.69%  "
$ws.Range("D12").Value = "'0.08010"
$ws.Range("E12").Value = "  -0.84%  "
$ws.Range("D13").Value = "'1.928.74"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").Value = "'5.346"
$ws.Range("E14").Value = "  -1.37%  "
$ws.Range("D15").Value = "'94.28"
$ws.Range("D16").Value = "'14.39"
$ws.Range("E16").Value = "  -2.89%  "
$ws.Range("D17").Value = "'30.337.50"
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("D18").Value = "'252.82"
$ws.Range("E18").Value = "  +0.36%  "
$ws.Range("D19").Value = "'0.000007924"
$ws.Range("E19").Value = "  -2.51%  "
$ws.Range("D20").Value = "'5.752"
$ws.Range("E20").Value = "  -0.56%  "
$ws.Range("D21").Value = "'2.188.13"
$ws.Range("D22").Value = "'0.9995"
$ws.Range("E22").Value = "  -0.11%  "
$ws.Range("D23").Value = "'1.0000"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").Value = "'6.669"
$ws.Range("E24").Value = "  -3.55%  "
$ws.Range("D25").Value = "'9.464"
$ws.Range("E25").Value = "  -2.57%  "
$ws.Range("D26").Value = "'165.51"
$ws.Range("E26").Value = "  +0.27%  "
$ws.Range("D27").Value = "'0.1332"
$ws.Range("E27").Value = "  +4.01%  "
$ws.Range("D28").Value = "'18.96"
$ws.Range("E28").Value = "  -1.38%  "
$ws.Range("D29").Value = "'2.179"
$ws.Range("E29").Value = "  -6.04%  "
$ws.Range("D30").Value = "'1.367"
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("D31").Value = "'1.518"
$ws.Range("E31").Value = "  -1.76%  "
$ws.Range("D32").Value = "'4.380"
$ws.Range("E32").Value = "  -1.14%  "
$ws.Range("D33").Value = "'4.119"
$ws.Range("E33").Value = "  -1.88%  "
$ws.Range("D34").Value = "'0.05152"
$ws.Range("E34").Value = "  -1.16%  "
$ws.Range("D35").Value = "'1.283"
$ws.Range("E35").Value = "  +1.06%  "
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("D37").Value = "'2.771"
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("D38").Value = "'0.01957"
$ws.Range("E38").Value = "  -0.07%  "
$ws.Range("D39").Value = "'2.803"
$ws.Range("E39").Value = "  +0.22%  "
$ws.Range("D40").Value = "'77.46"
$ws.Range("E40").Value = "  -2.06%  "
$ws.Range("D41").Value = "'6.407"
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("D42").Value = "'0.4453"
$ws.Range("E42").Value = "  -1.54%  "
$ws.Range("D43").Value = "'1.961"
$ws.Range("E43").Value = "  -3.31%  "
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("D45").Value = "'0.8333"
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("D46").Value = "'100.62"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("D47").Value = "'9.755"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").Value = "'7.463"
$ws.Range("E48").Value = "  +0.79%  "
$ws.Range("B49").Value = "Elrond"
$ws.Range("C49").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D49").Value = "'37.41"
$ws.Range("E49").Value = "  +2.40%  "
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "'980.56"
$ws.Range("E50").Value = "  +11.10%  "
$ws.Range("D51").Value = "'0.06014"
$ws.Range("E51").Value = "  -1.16%  "
